# Re-curate dimension columns that are actually measures: sector-descripcion,
# balance, ratios, rama-descripcion and pyg move from "iaest-dimension:*" /
# "dim" / "skos:Concept" (with an external mapping-*.xlsx) to
# "iaest-measure:*" / "medida" / "xsd:int" (no mapping file needed anymore).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(3, 11, 46, 64, 66)   # C, K, AT, BL, BN

foreach ($col in $cols) {
    # Row 2: iaest-dimension:<x> -> iaest-measure:<x>
    $cell2 = $ws.Cells.Item(2, $col)
    $cell2.Value = $cell2.Text -replace '^iaest-dimension:', 'iaest-measure:'

    # Row 3: dim -> medida
    $ws.Cells.Item(3, $col).Value = "medida"

    # Row 4: skos:Concept -> xsd:int
    $ws.Cells.Item(4, $col).Value = "xsd:int"

    # Row 5: drop the mapping-*.xlsx reference entirely (no longer needed)
    $ws.Cells.Item(5, $col).Clear()
}
